{"js": "// Replace the date line and the 25 three-digit \u00f7 one-digit problems with\n// the new day's values. Each old value is unique in the document, so a\n// direct search/replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-02-06 Tuesday\", \"2024-02-07 Wednesday\"],\n  [\"738\u00f78=\", \"760\u00f76=\"],\n  [\"879\u00f77=\", \"750\u00f77=\"],\n  [\"291\u00f79=\", \"922\u00f77=\"],\n  [\"111\u00f75=\", \"970\u00f72=\"],\n  [\"896\u00f78=\", \"331\u00f79=\"],\n  [\"605\u00f78=\", \"893\u00f78=\"],\n  [\"503\u00f75=\", \"157\u00f78=\"],\n  [\"408\u00f77=\", \"467\u00f74=\"],\n  [\"514\u00f78=\", \"239\u00f78=\"],\n  [\"362\u00f78=\", \"890\u00f76=\"],\n  [\"674\u00f72=\", \"321\u00f72=\"],\n  [\"314\u00f78=\", \"374\u00f76=\"],\n  [\"571\u00f79=\", \"338\u00f73=\"],\n  [\"524\u00f76=\", \"224\u00f72=\"],\n  [\"719\u00f73=\", \"782\u00f77=\"],\n  [\"142\u00f74=\", \"223\u00f74=\"],\n  [\"989\u00f78=\", \"864\u00f77=\"],\n  [\"658\u00f75=\", \"101\u00f73=\"],\n  [\"661\u00f77=\", \"996\u00f73=\"],\n  [\"329\u00f79=\", \"380\u00f79=\"],\n  [\"244\u00f77=\", \"976\u00f75=\"],\n  [\"743\u00f74=\", \"401\u00f77=\"],\n  [\"399\u00f79=\", \"606\u00f75=\"],\n  [\"383\u00f73=\", \"934\u00f77=\"],\n  [\"654\u00f74=\", \"120\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 three-digit \u00f7 one-digit problems with\n# the new day's values. Each old value is unique in the document, so a\n# direct Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-06 Tuesday\", \"2024-02-07 Wednesday\"),\n    @(\"738\u00f78=\", \"760\u00f76=\"),\n    @(\"879\u00f77=\", \"750\u00f77=\"),\n    @(\"291\u00f79=\", \"922\u00f77=\"),\n    @(\"111\u00f75=\", \"970\u00f72=\"),\n    @(\"896\u00f78=\", \"331\u00f79=\"),\n    @(\"605\u00f78=\", \"893\u00f78=\"),\n    @(\"503\u00f75=\", \"157\u00f78=\"),\n    @(\"408\u00f77=\", \"467\u00f74=\"),\n    @(\"514\u00f78=\", \"239\u00f78=\"),\n    @(\"362\u00f78=\", \"890\u00f76=\"),\n    @(\"674\u00f72=\", \"321\u00f72=\"),\n    @(\"314\u00f78=\", \"374\u00f76=\"),\n    @(\"571\u00f79=\", \"338\u00f73=\"),\n    @(\"524\u00f76=\", \"224\u00f72=\"),\n    @(\"719\u00f73=\", \"782\u00f77=\"),\n    @(\"142\u00f74=\", \"223\u00f74=\"),\n    @(\"989\u00f78=\", \"864\u00f77=\"),\n    @(\"658\u00f75=\", \"101\u00f73=\"),\n    @(\"661\u00f77=\", \"996\u00f73=\"),\n    @(\"329\u00f79=\", \"380\u00f79=\"),\n    @(\"244\u00f77=\", \"976\u00f75=\"),\n    @(\"743\u00f74=\", \"401\u00f77=\"),\n    @(\"399\u00f79=\", \"606\u00f75=\"),\n    @(\"383\u00f73=\", \"934\u00f77=\"),\n    @(\"654\u00f74=\", \"120\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format, ReplaceWith,\n    # Replace(2=wdReplaceAll)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
